# Adding TWB Primary Nominated Professional Consent Date.
#
# This script reproduces the data changes: a set of "Collection Occasion"
# keys (CO04..CO13, CO17) used by the measure sheets are renamed to their
# "-1" suffixed counterparts (CO04-1..CO13-1), a couple of header / key
# fixes are applied (TWB Recommendation Outs header, a new TWB Plan
# TWBP04), and eight new Collection Occasion rows are appended to the
# "Collection Occasions" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# K10+ : collection_occasion references CO08..CO11 -> CO08-1..CO11-1
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("K10+")
$s.Range("C2").Value = "CO08-1"
$s.Range("C3").Value = "CO09-1"
$s.Range("C4").Value = "CO10-1"
$s.Range("C5").Value = "CO11-1"
$s.Range("A2:R5").Select() | Out-Null

# ---------------------------------------------------------------------
# K5 : collection_occasion references CO04, CO05, CO12, CO13 -> *-1
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("K5")
$s.Range("C2").Value = "CO04-1"
$s.Range("C3").Value = "CO05-1"
$s.Range("C4").Value = "CO12-1"
$s.Range("C5").Value = "CO13-1"
$s.Range("A1:J5").Select() | Out-Null

# ---------------------------------------------------------------------
# SDQ : collection_occasion references CO07, CO17 -> CO07-1, CO13-1
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("SDQ")
$s.Range("C3").Value = "CO07-1"
$s.Range("C4").Value = "CO13-1"
$s.Range("A2:BB4").Select() | Out-Null

# ---------------------------------------------------------------------
# WHO-5 : collection_occasion references CO08..CO11 -> CO08-1..CO11-1
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("WHO-5")
$s.Range("C2").Value = "CO08-1"
$s.Range("C3").Value = "CO09-1"
$s.Range("C4").Value = "CO10-1"
$s.Range("C5").Value = "CO11-1"
$s.Range("A2:J5").Select() | Out-Null

# ---------------------------------------------------------------------
# SIDAS : no data changes, selection only
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("SIDAS")
$s.Range("A2:I3").Select() | Out-Null

# ---------------------------------------------------------------------
# TWB Plans : collection_occasion references CO08..CO11 -> CO08-1..CO11-1
# plus a new plan key TWBP04 replacing TWBP03 on row 5 (and row4's
# duplicate TWBP02 becomes TWBP03)
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("TWB Plans")
$s.Range("C2").Value = "CO08-1"
$s.Range("C3").Value = "CO09-1"
$s.Range("B4").Value = "TWBP03"
$s.Range("C4").Value = "CO10-1"
$s.Range("B5").Value = "TWBP04"
$s.Range("C5").Value = "CO11-1"
$s.Range("A2:E5").Select() | Out-Null

# ---------------------------------------------------------------------
# TWB NIs : collection_occasion references CO08..CO11 -> CO08-1..CO11-1
# This sheet ends up being the active tab.
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("TWB NIs")
$s.Range("C2").Value = "CO08-1"
$s.Range("C3").Value = "CO09-1"
$s.Range("C4").Value = "CO10-1"
$s.Range("C5").Value = "CO11-1"

# ---------------------------------------------------------------------
# TWB Recommendation Outs : fix header B1 (was mislabeled with the
# "TWB Critical Incidents" key name)
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("TWB Recommendation Outs")
$s.Range("B1").Value = "twb_recommendation_out_key"
$s.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------
# Collection Occasions : append eight new rows (the new TWB Primary
# Nominated Professional Consent Date collection occasions CO06-1..CO13-1)
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("Collection Occasions")

$s.Cells.Item(4,1).Value = "PHN999:NFP01"
$s.Cells.Item(4,2).Value = "CO06-1"
$s.Cells.Item(4,3).Value = "CL0001-E01"
$s.Cells.Item(4,4).Value = 16062018
$s.Cells.Item(4,5).Value = 2
$s.Cells.Item(4,6).Value = "tag1 tag2"

$s.Cells.Item(5,1).Value = "PHN999:NFP01"
$s.Cells.Item(5,2).Value = "CO07-1"
$s.Cells.Item(5,3).Value = "CL0001-E01"
$s.Cells.Item(5,4).Value = 16062019
$s.Cells.Item(5,5).Value = 2
$s.Cells.Item(5,6).Value = "tag1 tag2"

$s.Cells.Item(6,1).Value = "PHN999:NFP01"
$s.Cells.Item(6,2).Value = "CO08-1"
$s.Cells.Item(6,3).Value = "CL0002-E01"
$s.Cells.Item(6,4).Value = 16062016
$s.Cells.Item(6,5).Value = 1
$s.Cells.Item(6,6).Value = "tag1 tag2"

$s.Cells.Item(7,1).Value = "PHN999:NFP01"
$s.Cells.Item(7,2).Value = "CO09-1"
$s.Cells.Item(7,3).Value = "CL0001-E01"
$s.Cells.Item(7,4).Value = 16062018
$s.Cells.Item(7,5).Value = 2
$s.Cells.Item(7,6).Value = "tag1 tag2"

$s.Cells.Item(8,1).Value = "PHN999:NFP01"
$s.Cells.Item(8,2).Value = "CO10-1"
$s.Cells.Item(8,3).Value = "CL0002-E01"
$s.Cells.Item(8,4).Value = 20062016
$s.Cells.Item(8,5).Value = 2
$s.Cells.Item(8,6).Value = "tag1 tag2"

$s.Cells.Item(9,1).Value = "PHN999:NFP01"
$s.Cells.Item(9,2).Value = "CO11-1"
$s.Cells.Item(9,3).Value = "CL0002-E01"
$s.Cells.Item(9,4).Value = 16062016
$s.Cells.Item(9,5).Value = 2
$s.Cells.Item(9,6).Value = "tag1 tag2"

$s.Cells.Item(10,1).Value = "PHN999:NFP01"
$s.Cells.Item(10,2).Value = "CO12-1"
$s.Cells.Item(10,3).Value = "CL0001-E01"
$s.Cells.Item(10,4).Value = 16062016
$s.Cells.Item(10,5).Value = 2
$s.Cells.Item(10,6).Value = "tag1 tag2"

$s.Cells.Item(11,1).Value = "PHN999:NFP01"
$s.Cells.Item(11,2).Value = "CO13-1"
$s.Cells.Item(11,3).Value = "CL0002-E01"
$s.Cells.Item(11,4).Value = 16062016
$s.Cells.Item(11,5).Value = 2
$s.Cells.Item(11,6).Value = "tag1 tag2"

$s.Range("A2:F11").Select() | Out-Null

# ---------------------------------------------------------------------
# Metadata : no longer the active tab
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("Metadata")
$s.Range("B4").Select() | Out-Null

# ---------------------------------------------------------------------
# Finally, make "TWB NIs" the active sheet/tab (matches the saved view
# state in the workbook: activeTab points at this sheet and it carries
# tabSelected="1").
# ---------------------------------------------------------------------
$s = $wb.Worksheets.Item("TWB NIs")
$s.Activate()
$s.Range("A2:E5").Select() | Out-Null
